# bugfix - extra page break
# fixed extra page break
#
# The "EXPERIENCE, continued (page 2)" paragraph on the resume starts
# with a stray manual page break, and it is preceded by two identical
# (duplicate) blank paragraphs instead of one. This produced an extra,
# unwanted blank page in the rendered document. Fix:
#   1. Delete the leading manual page-break run in the
#      "EXPERIENCE, continued (page 2)" paragraph.
#   2. Delete one of the two duplicate blank paragraphs that precede it.

$d = $word.ActiveDocument

# --- Locate the "EXPERIENCE, continued (page 2)" paragraph ---
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*EXPERIENCE, continued (page 2)*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq 0) {
    throw "Could not find 'EXPERIENCE, continued (page 2)' paragraph"
}

# --- Remove the leading manual page-break character (chr 12 / ^12) ---
$targetPara = $d.Paragraphs.Item($targetIndex)
$breakRange = $targetPara.Range.Duplicate
$breakRange.SetRange($targetPara.Range.Start, $targetPara.Range.Start + 1)
if ($breakRange.Text -eq [char]12) {
    $breakRange.Delete()
}

# --- Remove the extra duplicate blank paragraph right before it ---
# Walking back from the target: [target-1] is a blank "heading spacer"
# paragraph, and [target-2]/[target-3] are two back-to-back blank
# paragraphs with identical formatting (the duplicate introduced by the
# bug). Delete the later (closer) one of that duplicate pair, [target-2],
# leaving a single blank paragraph behind -- matching how the fix shipped.
$dupIndex = $targetIndex - 2
$priorIndex = $dupIndex - 1

$dupPara = $d.Paragraphs.Item($dupIndex)
$priorPara = $d.Paragraphs.Item($priorIndex)

$dupIsBlank = ($dupPara.Range.Text.Trim().Length -eq 0)
$priorIsBlank = ($priorPara.Range.Text.Trim().Length -eq 0)

if ($dupIsBlank -and $priorIsBlank) {
    $dupPara.Range.Delete()
} elseif ($dupIsBlank) {
    # Fallback: structure differs from the expected duplicate-pair shape,
    # but the immediate predecessor is still blank -- remove it, since
    # that is the paragraph the page-break bugfix targets.
    $dupPara.Range.Delete()
}
